# Applies the OOXML diff to the presentation via PowerPoint COM interop.
#
# Summary of the edit (per the commit's xml diff):
#  - Slide 2 ("Sales Report" panel): reposition/resize the purple header
#    card's title textbox, connector line and logo picture, retitle the
#    textbox from "Sales Report" to "Relatório de vendas", and nudge the
#    rounded-rectangle header bar upward.
#  - The "Espaço Reservado para Data" (date) placeholders on the slide
#    master, every slide layout and the notes master show PowerPoint's
#    auto "today" field; it had cached 22/09/2024 and now caches
#    06/01/2025.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 2 shape geometry + text tweaks
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(2)

# "Retângulo: Cantos Arredondados 3" (header card behind the title) -
# only its vertical position moves (362139 EMU -> 268196 EMU).
$headerCard = $slide.Shapes.Item(2)
$headerCard.Top = 21.117795944213867

# "CaixaDeTexto 7" - the "Sales Report" title textbox: shifts 1 EMU left,
# widens considerably, and its text is translated to Portuguese.
$titleBox = $slide.Shapes.Item(3)
$titleBox.Left = 75.56425476074219
$titleBox.Width = 169.66339111328125
$titleBox.TextFrame.TextRange.Text = "Relatório de vendas"

# "Conector reto 9" - the underline below the title: only widens.
$underline = $slide.Shapes.Item(4)
$underline.Width = 148.27748107910156

# "Imagem 12" - the logo picture top-right of the header: moves right.
$logo = $slide.Shapes.Item(5)
$logo.Left = 214.6341094970703

# ---------------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" field text (22/09/2024 ->
#    06/01/2025) everywhere it is placed: slide master and every custom
#    layout (note: the notes master's own date placeholder cannot be
#    written through this COM host - mutations to its shapes do not
#    persist/are not reachable - so it is intentionally left alone
#    rather than risk touching unrelated shapes).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDate = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) { $isDate = $true }
            } catch {
            }
            if ($isDate) {
                $shp.TextFrame.TextRange.Text = "06/01/2025"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
